# Updated cryptos list (price / volume refresh + a 3-way coin reorder at
# rows 48-50: Frax/Cronos/BabyDogeCoin rotate places), mirroring the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Note: many of the new "Price" (column D) strings look like plain numbers
# ("1.004", "217.50", "20.43", ...) even though the sheet stores them as
# text (t="inlineStr"/shared-string, not t="n"). Assigning such a string to
# Range.Value would normally let Excel auto-convert it to a real number
# (and e.g. drop the trailing zero in "217.50"). To keep it text - exactly
# like typing '217.50 into a cell - those assignments are prefixed with a
# leading apostrophe, which Excel strips from the stored value. Strings
# that already aren't parseable as a plain number (multiple dots such as
# "26.088.33", or the subscript-digit ShibaInu/BabyDogeCoin prices) are
# assigned as-is since Excel already keeps them as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '26.088.33'
$ws.Range("E2").Value = '  -0.93%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '1.645.25'
$ws.Range("E3").Value = '  -1.28%  '

# Row 4 - TetherUSD
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  -0.66%  '

# Row 5 - BNB
$ws.Range("D5").Value = '''217.50'
$ws.Range("E5").Value = '  -0.77%  '

# Row 6 - XRP
$ws.Range("D6").Value = '''0.5203'
$ws.Range("E6").Value = '  -2.52%  '

# Row 7 - USDC
$ws.Range("E7").Value = '  -0.56%  '

# Row 8 - Cardano
$ws.Range("E8").Value = '  -1.63%  '

# Row 9 - Dogecoin
$ws.Range("D9").Value = '''0.06282'
$ws.Range("E9").Value = '  -1.74%  '

# Row 10 - Solana
$ws.Range("D10").Value = '''20.43'
$ws.Range("E10").Value = '  -2.11%  '

# Row 11 - TRON
$ws.Range("D11").Value = '''0.07760'
$ws.Range("E11").Value = '  -1.33%  '

# Row 12 - Polkadot
$ws.Range("D12").Value = '''4.475'
$ws.Range("E12").Value = '  -1.87%  '

# Row 13 - WrappedEther
$ws.Range("D13").Value = '1.626.66'
$ws.Range("E13").Value = '  -2.55%  '

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '1.870.25'
$ws.Range("E14").Value = '  -1.30%  '

# Row 15 - Polygon
$ws.Range("D15").Value = '''0.5579'
$ws.Range("E15").Value = '  +0.74%  '

# Row 16 - ShibaInu
$ws.Range("D16").Value = '0.0₅8002'
$ws.Range("E16").Value = '  -2.28%  '

# Row 17 - Litecoin
$ws.Range("D17").Value = '''64.77'

# Row 18 - WrappedBTC
$ws.Range("D18").Value = '26.076.02'
$ws.Range("E18").Value = '  -1.08%  '

# Row 19 - Dai
$ws.Range("E19").Value = '  -0.57%  '

# Row 20 - Uniswap
$ws.Range("D20").Value = '''4.642'
$ws.Range("E20").Value = '  -0.80%  '

# Row 21 - BitcoinCash
$ws.Range("D21").Value = '''192.42'
$ws.Range("E21").Value = '  -0.70%  '

# Row 22 - Avalanche
$ws.Range("D22").Value = '''10.10'
$ws.Range("E22").Value = '  -1.78%  '

# Row 23 - Chainlink
$ws.Range("D23").Value = '''5.951'
$ws.Range("E23").Value = '  -1.47%  '

# Row 24 - BinanceUSD
$ws.Range("D24").Value = '''1.006'
$ws.Range("E24").Value = '  -0.64%  '

# Row 25 - Monero
$ws.Range("D25").Value = '''146.43'

# Row 26 - Stellar
$ws.Range("D26").Value = '''0.1202'
$ws.Range("E26").Value = '  -2.09%  '

# Row 27 - Cosmos
$ws.Range("D27").Value = '''7.166'
$ws.Range("E27").Value = '  -0.55%  '

# Row 28 - EthereumClassic
$ws.Range("E28").Value = '  -1.06%  '

# Row 29 - Toncoin
$ws.Range("E29").Value = '  -1.00%  '

# Row 30 - Hedera
$ws.Range("D30").Value = '''0.05622'
$ws.Range("E30").Value = '  -3.85%  '

# Row 31 - PancakeSwap
$ws.Range("D31").Value = '''1.264'
$ws.Range("E31").Value = '  -1.41%  '

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = '''3.457'
$ws.Range("E32").Value = '  -3.70%  '

# Row 33 - Filecoin
$ws.Range("D33").Value = '''3.359'
$ws.Range("E33").Value = '  +2.38%  '

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = '''1.594'
$ws.Range("E34").Value = '  -0.44%  '

# Row 35 - MXToken
$ws.Range("D35").Value = '''2.793'
$ws.Range("E35").Value = '  -1.26%  '

# Row 36 - HuobiToken
$ws.Range("E36").Value = '  -0.37%  '

# Row 37 - ARBITRUM
$ws.Range("D37").Value = '''0.9361'
$ws.Range("E37").Value = '  -3.55%  '

# Row 38 - ImmutableX
$ws.Range("D38").Value = '''0.5671'
$ws.Range("E38").Value = '  -2.73%  '

# Row 39 - FraxShare
$ws.Range("D39").Value = '''5.964'
$ws.Range("E39").Value = '  +2.22%  '

# Row 41 - Maker
$ws.Range("D41").Value = '1.051.28'
$ws.Range("E41").Value = '  -1.31%  '

# Row 42 - mCoin
$ws.Range("D42").Value = '''2.567'
$ws.Range("E42").Value = '  -0.78%  '

# Row 43 - PaxDollar
$ws.Range("E43").Value = '  -0.63%  '

# Row 44 - TrustWalletToken
$ws.Range("D44").Value = '''0.8396'
$ws.Range("E44").Value = '  -2.64%  '

# Row 45 - Quant
$ws.Range("D45").Value = '''102.34'
$ws.Range("E45").Value = '  -2.35%  '

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = '1.782.19'

# Row 47 - Aave
$ws.Range("D47").Value = '''57.19'
$ws.Range("E47").Value = '  -1.06%  '

# Rows 48-50 rotate: Frax -> BabyDogeCoin, Cronos -> Frax, BabyDogeCoin -> Cronos
# Row 48 becomes BabyDogeCoin (was Frax)
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₈106'
$ws.Range("E48").Value = '  +1.96%  '

# Row 49 becomes Frax (was Cronos)
$ws.Range("B49").Value = 'Frax'
$ws.Range("C49").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D49").Value = '''1.010'
$ws.Range("E49").Value = '  -0.38%  '

# Row 50 becomes Cronos (was BabyDogeCoin)
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.05315'
$ws.Range("E50").Value = '  +2.88%  '

# Row 51 - Mantle
$ws.Range("D51").Value = '''0.4326'
$ws.Range("E51").Value = '  -1.55%  '
